# Weekly update: insert a new data row at the top of the "Zapallo italiano"
# price series (row 533), pushing the existing rows 533:615 down to 534:616.
# The new row carries the latest weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 533 - everything below shifts down
# by one (old row 533 becomes 534, ..., old row 615 becomes 616).
$ws.Rows("533:533").Insert()

# Populate the newly inserted row 533 with the new weekly record.
$ws.Range("A533").Value = 9
$ws.Range("B533").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C533").Value = "Metropolitana"
$ws.Range("D533").Value = 45180
$ws.Range("E533").Value = 13
$ws.Range("F533").Value = 100112032
$ws.Range("G533").Value = "Zapallo italiano"
$ws.Range("H533").Value = "Sin especificar"
$ws.Range("I533").Value = "Primera"
$ws.Range("J533").Value = 70
$ws.Range("K533").Value = 15000
$ws.Range("L533").Value = 16000
$ws.Range("M533").Value = 15500
$ws.Range("N533").Value = "`$/caja 50 unidades"
$ws.Range("O533").Value = "Región de Arica y Parinacota"
$ws.Range("P533").Value = 310
$ws.Range("Q533").Value = 50
$ws.Range("R533").Value = "Hortaliza"
